# Update báo cáo tiến độ lần 1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project 1")

# Row 9: mark the task 100% complete and add a hyperlink with a new URL in column G
$ws.Range("B9").Value = 1
$ws.Range("G9").Hyperlinks.Add($ws.Range("G9"), "http://localhost:52360/api/luatxettuyen/create?sup=30&&con=50") | Out-Null

# Row 10 & 11: mark tasks 80% complete
$ws.Range("B10").Value = 0.8
$ws.Range("B11").Value = 0.8

# Row 13: mark task 100% complete
$ws.Range("B13").Value = 1

# Update selection / view position to match the saved state
$ws.Activate()
$ws.Range("G10").Select()
$ws.Application.ActiveWindow.ScrollRow = 7
